# Events.pptx edit: add plain-language descriptions, simplify level text,
# enlarge / reflow the condition+outcome boxes and append a "Hintergrund"
# (background) textbox to each of the five event-card slides.
#
# Slide map in this deck (slides are 1-indexed):
#   2 -> OEM-Audit               (has BOTH Erfuellung + Nicht-Erfuellung outcomes)
#   3 -> Personalwechsel         (only Nicht-Erfuellung outcome, shown on the right)
#   4 -> DSGVO-Bonus             (only Erfuellung outcome, shown on the left)
#   5 -> Investoren-Vertrauen    (only Erfuellung outcome, shown on the left)
#   6 -> Compliance-Luecke       (only Nicht-Erfuellung outcome, shown on the right)

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Convert an OOXML EMU value to the Points unit the COM object model uses for
# Left/Top/Width/Height. A tiny epsilon compensates for the single-precision
# float round-trip the host performs internally (it truncates), so values
# that are exact multiples of 12700 EMU come back out exact too.
function EMU($emu) {
    return ($emu / 12700.0) + 0.00005
}

# Reposition/resize a shape precisely, in EMU.
function Set-ShapeRect($shape, $x, $y, $cx, $cy) {
    $shape.Left = EMU($x)
    $shape.Top = EMU($y)
    $shape.Width = EMU($cx)
    $shape.Height = EMU($cy)
}

# RGB hex string ("RRGGBB") -> the BGR-packed integer PowerPoint's
# Font.Color.RGB / Fill.ForeColor.RGB setters expect.
function RGBHex($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Per-slide text content (old condition/description text is already in the
# deck and only needs new sizes; only the strings below are new).
# ---------------------------------------------------------------------------

$slideData = @(
    @{
        Index    = 2
        DescNew  = "Der große Automobilkunde prüft die Sicherheitsstandards. Wer gut vorbereitet ist, gewinnt Vertrauen - wer schlecht abschneidet, riskiert Aufträge."
        CondNew  = "Wurde der Sicherheits-Zielwert nach Welle 1 erreicht?"
        HintNew  = "Hintergrund: Bestanden: Kundenvertrauen steigt deutlich. Nicht bestanden: Vertrauen sinkt, Folgeaufträge gefährdet."
        BothOutcomes = $true
    },
    @{
        Index    = 3
        DescNew  = "Ohne regelmäßige Schulungen und Sensibilisierung werden IT-Sicherheitsaufgaben zur Belastung. Überlastete Mitarbeiter verlassen das Unternehmen."
        CondNew  = "Wurden die Mitarbeiter nicht ausreichend geschult (Awareness unter Level 2)?"
        HintNew  = "Hintergrund: Fluktuation kostet Geld (Einarbeitung) und Know-how geht verloren."
        BothOutcomes = $false
    },
    @{
        Index    = 4
        DescNew  = "Gute Zugriffskontrolle und Logging sind die Basis für Datenschutz-Compliance. Wer hier investiert hat, wird belohnt."
        CondNew  = "Sind Zugriffskontrolle (M1) UND Logging (M2) mindestens auf Level 2?"
        HintNew  = "Hintergrund: Compliance-Nachweis erleichtert Kundengewinnung und spart Bußgelder."
        BothOutcomes = $false
    },
    @{
        Index    = 5
        DescNew  = "Wer viel in Sicherheit investiert, zeigt Weitsicht. Das überzeugt Investoren und Geldgeber."
        CondNew  = "Wurde das höchste Budget-Level gewählt?"
        HintNew  = "Hintergrund: Hohe Sicherheitsinvestitionen signalisieren professionelles Risikomanagement."
        BothOutcomes = $false
    },
    @{
        Index    = 6
        DescNew  = "Zu wenig Budget bedeutet Kompromisse bei der Sicherheit. Das fällt spätestens bei Audits und Kundenanfragen negativ auf."
        CondNew  = "Wurde das niedrigste Budget-Level gewählt?"
        HintNew  = "Hintergrund: Sichtbare Sicherheitslücken schaden dem Ruf bei Kunden und Partnern."
        BothOutcomes = $false
    }
)

foreach ($cfg in $slideData) {

    $s = $p.Slides.Item($cfg.Index)

    # --- Header bar (orange rectangle) shrinks a bit ----------------------
    $rect1 = $s.Shapes.Item("Rectangle 1")
    Set-ShapeRect $rect1 0 0 9144000 1371600

    # --- Title ("TextBox 2") moves up / shrinks ----------------------------
    $title = $s.Shapes.Item("TextBox 2")
    Set-ShapeRect $title 457200 182880 8229600 548640
    $title.TextFrame.TextRange.Font.Size = 36

    # --- Trigger line ("TextBox 3") moves up, smaller font -----------------
    $trigger = $s.Shapes.Item("TextBox 3")
    Set-ShapeRect $trigger 457200 777240 8229600 457200
    $trigger.TextFrame.TextRange.Font.Size = 18

    # --- Plain-language description ("TextBox 4") --------------------------
    $desc = $s.Shapes.Item("TextBox 4")
    $desc.TextFrame.TextRange.Text = $cfg.DescNew
    $desc.TextFrame.TextRange.Font.Size = 18
    Set-ShapeRect $desc 457200 1554480 8229600 822960

    # --- "Bedingung" card background ("Rounded Rectangle 5") ---------------
    $condBg = $s.Shapes.Item("Rounded Rectangle 5")
    Set-ShapeRect $condBg 457200 2514600 8229600 1188720

    # --- "Bedingung" text ("TextBox 6": 2 paragraphs) -----------------------
    $condBox = $s.Shapes.Item("TextBox 6")
    $condTr = $condBox.TextFrame.TextRange
    $condPara1 = $condTr.Paragraphs(1, 1)
    $condPara1.Font.Size = 14
    $condPara2 = $condTr.Paragraphs(2, 1)
    $condPara2.Text = $cfg.CondNew
    $condPara2.Font.Size = 16
    Set-ShapeRect $condBox 640080 2606040 7863840 1005840

    if ($cfg.BothOutcomes) {
        # Slide has both a green "Bei Erfuellung" box (left) and a red
        # "Bei Nicht-Erfuellung" box (right); neither moves horizontally.

        $okBg = $s.Shapes.Item("Rounded Rectangle 7")
        Set-ShapeRect $okBg 457200 3840480 3931920 1280160

        $okBox = $s.Shapes.Item("TextBox 8")
        $okTr = $okBox.TextFrame.TextRange
        $okTr.Paragraphs(1, 1).Font.Size = 13
        $okTr.Paragraphs(2, 1).Font.Size = 18
        Set-ShapeRect $okBox 594360 3931920 3657600 1097280

        $badBg = $s.Shapes.Item("Rounded Rectangle 9")
        Set-ShapeRect $badBg 4754880 3840480 3931920 1280160

        $badBox = $s.Shapes.Item("TextBox 10")
        $badTr = $badBox.TextFrame.TextRange
        $badTr.Paragraphs(1, 1).Font.Size = 13
        $badTr.Paragraphs(2, 1).Font.Size = 18
        Set-ShapeRect $badBox 4892040 3931920 3657600 1097280
    }
    else {
        # Slide has a single outcome box (either green or red); it widens
        # and re-centers.
        $outBg = $s.Shapes.Item("Rounded Rectangle 7")
        Set-ShapeRect $outBg 2286000 3840480 4572000 1280160

        $outBox = $s.Shapes.Item("TextBox 8")
        $outTr = $outBox.TextFrame.TextRange
        $outTr.Paragraphs(1, 1).Font.Size = 13
        $outTr.Paragraphs(2, 1).Font.Size = 18
        Set-ShapeRect $outBox 2423160 3931920 4297680 1097280
    }

    # --- New "Hintergrund" textbox at the bottom of the card ---------------
    $hint = $s.Shapes.AddTextbox(1, EMU(457200), EMU(5303520), EMU(8229600), EMU(731520))
    $hint.Fill.Visible = $false
    $hint.TextFrame.WordWrap = -1
    $hint.TextFrame.AutoSize = 1
    $hintTr = $hint.TextFrame.TextRange
    $hintTr.Text = $cfg.HintNew
    $hintTr.ParagraphFormat.Alignment = 2
    $hintTr.Font.Size = 14
    $hintTr.Font.Italic = $true
    $hintTr.Font.Color.RGB = RGBHex("343A40")
    Set-ShapeRect $hint 457200 5303520 8229600 731520
}
